$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('Pick-A-Stone', '5701 Bingle Rd Ste A106, Houston, TX 77092'),
  @('Absolute Marble and Granite, Inc.', '505 Parker St, Houston, TX 77007'),
  @('Earth Stone Colours', '8758 Clay Rd Suite #400, Houston, TX 77080'),
  @('All World Stone Supply', '23527 FM1314 Suite B, Porter, TX 77365'),
  @('Instyle Granite & Marble', '14210 Aston St, Houston, TX 77040'),
  @('Marvel-Works.Com', '6700 Woodlands Pkwy STE 230 #455, The Woodlands, TX 77382'),
  @('Earth Stone Colours', '4305 B, 4305 Richmond Ave. Ste A, Houston, TX 77027'),
  @('Future Stone', '6040 Brittmoore Rd, Houston, TX 77041'),
  @('ABC STONE & KITCHEN', '335 Staffordshire Rd STE. 5, Stafford, TX 77477'),
  @('Earth Stone Colours (Prefab)', '8700 Clay Rd Suite #230, Houston, TX 77080'),
  @('Marble Line', '32522 McAllister Rd, Brookshire, TX 77423'),
  @('Granite Countertops Houston Tx', 'granitecountertopshoustontx.com'),
  @('Millstone Marble', '20737 FM 1774, Plantersville, TX 77363'),
  @('Mega Stone Inc', 'megastoneinc.business.site'),
  @('Medina Countertops', 'medinacountertops.com'),
  @('TT Stone Outlet', '5750 Campbell Rd, Houston, TX 77041'),
  @('JED Custom Countertops Corp.', '27405 Old Washington Rd, Hockley, TX 77447'),
  @('TaylorMade Custom Granite', '714 Honea Egypt Rd suite f, Magnolia, TX 77354'),
  @('granite countertops starting at $5.99 per sf granite and cabinet warehouse', '8307 Bauman Rd, Houston, TX 77022'),
  @('Xiamen Superior Stone Co.,Ltd.', '12333 Sowden Rd Ste B # 30327, Houston, TX 77080')
)

$startRow = 98
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

Write-Host "Done. Rows written:" $data.Count
